# Insert a new record row at row 47 (pushes existing rows 47-64 down to 48-65)
# and populate it with the new "Espárragos" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(47, 3).Value = "Maule"
$ws.Cells.Item(47, 4).Value = 44830
$ws.Cells.Item(47, 5).Value = 7
$ws.Cells.Item(47, 6).Value = 300000000
$ws.Cells.Item(47, 7).Value = "Espárragos"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 500
$ws.Cells.Item(47, 11).Value = 2000
$ws.Cells.Item(47, 12).Value = 2000
$ws.Cells.Item(47, 13).Value = 2000
$ws.Cells.Item(47, 14).Value = "$/kilo"
$ws.Cells.Item(47, 15).Value = "Provincia de Linares"
$ws.Cells.Item(47, 16).Value = 2000
$ws.Cells.Item(47, 17).Value = 1
$ws.Cells.Item(47, 18).Value = "Hortaliza"
